# Update review sheets for CYRS and HSI
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Introduction "
$ws2 = $wb.Worksheets.Item(2)   # "Cross review points "

# --- sheet1 "Introduction " edits ---

# Document status table: "Ref Version" bumped from V1.5 to V1.6
$ws1.Range("D7").Value = "V1.6"

# Document status table: "Last update" changed from a date serial to literal text "14/2/2020"
$ws1.Range("D9").Value = "14/2/2020"

# New row 15 entry added to the Document History table
$ws1.Range("B15").Value = 0.3
$ws1.Range("C15").Value = "T.Sharaby"
$ws1.Range("E15").Value = "14/2/202"
$ws1.Range("G15").Value = "Update the status of the last open point"

# --- sheet2 "Cross review points " edits ---

# Row 9: Status changed from "Open" to "Resolved"
$ws2.Range("H9").Value = "Resolved"

# Row 9: Comment text updated with TSH follow-up note
$ws2.Range("I9").Value = "Added the needed requirements`nTSH: I mean here the system doesn't descrip those use cases , are you able to handle all these use cases in the SRS req , if yes then it is ok for me ?"

# Row 9: Comment cell alignment switched from centered to left-aligned
$ws2.Range("I9").HorizontalAlignment = -4131

# --- Sheet view / selection updates ---

# "Cross review points " is no longer the active tab; select cell H12 on it
[void]$ws2.Activate()
[void]$ws2.Range("H12").Select()

# "Introduction " becomes the active tab, with B10:H10 selected
[void]$ws1.Activate()
[void]$ws1.Range("B10:H10").Select()
